$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look numeric must stay as text (matching original inlineStr cells)
# and keep the default (unstyled) cell formatting, so we temporarily mark them as Text,
# assign the value, then restore the "Normal" style so no stray style index remains.
$textForceCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D17", "D18", "D19", "D23", "D27", "D28", "D29", "D30", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates (new values taken from the target revision)
$ws.Range("D2").Value = "28.909.92"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.887.16"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "331.63"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D7").Value = "0.4616"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").Value = "0.4096"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("D9").Value = "47.48"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "0.07986"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "0.9906"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "21.71"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "1.874.61"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "7.068"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "89.14"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "0.00001029"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "0.06563"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "28.943.97"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").Value = "5.375"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "2.095.21"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "157.59"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "19.67"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "2.118"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "5.408"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "0.9770"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "0.09354"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("E34").Value = "  +3.06%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "5.278"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "0.06058"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "0.02229"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "8.253"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "1.175"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "0.5771"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "10.13"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "1.268"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "2.281"
$ws.Range("E46").Value = "  +10.56%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5473"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "11.95"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "1.906"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "0.07013"
$ws.Range("E50").Value = "  -7.12%  "
$ws.Range("D51").Value = "110.62"
$ws.Range("E51").Value = "  -1.08%  "

# Restore default styling on the cells that were temporarily marked as Text,
# so their style index matches the original (unstyled) cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
